$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: (row, E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T)
# Columns E..T correspond to spreadsheet columns 5..20.
$data = @(
  @(2, 2, 1, 38.38922700000001, 76.77845400000001, 0.4452295445267456, 0.3751312191747254, 2, 1, 1.0305975, 2.061195, 0.09953130389913815, 0.07508903821244231, 39.5638413731325, 158.25536549253, 0.04431427710116637, 0.02816824245129103),
  @(3, 2, 1, 38.38922700000001, 76.77845400000001, 0.4452295445267456, 0.3751312191747254, 3, 1, 1.525218666666667, 4.575656, 0.147299991145562, 0.1666904917928634, 58.55196561930401, 351.3117937158241, 0.06558230796653221, 0.06253080741109139),
  @(4, 2, 1, 38.38922700000001, 76.77845400000001, 0.4452295445267456, 0.3751312191747254, 3, 1, 1.729477333333333, 5.188432, 0.1670265395080728, 0.189013833582295, 66.39329794068802, 398.359787644128, 0.07436515010905771, 0.07090498983261496),
  @(5, 2, 1, 38.38922700000001, 76.77845400000001, 0.4452295445267456, 0.3751312191747254, 3, 1, 2.911407000000001, 8.734221000000002, 0.2811729456854671, 0.3181864182791614, 111.766664212389, 670.5999852743342, 0.1251865025407839, 0.119361659013901),
  @(6, 2, 1, 38.38922700000001, 76.77845400000001, 0.4452295445267456, 0.3751312191747254, 3, 1, 0.5748966666666666, 1.72469, 0.05552139883960665, 0.06283020932741304, 22.06983863821, 132.41903182926, 0.02471976711684585, 0.02356957302599566),
  @(7, 2, 1, 38.38922700000001, 76.77845400000001, 0.4452295445267456, 0.3751312191747254, 2, 1, 2.582909, 5.165818, 0.2494478209221534, 0.1881900088058249, 99.15587992134301, 396.6235196853721, 0.1110615396923596, 0.07059594743983141),
  @(8, 3, 1, 3.154770000000001, 9.464310000000001, 0.03658830666704076, 0.04624159466596638, 2, 1, 1.0305975, 2.061195, 0.09953130389913815, 0.07508903821244231, 3.251298075075001, 19.50778845045, 0.003641681870032096, 0.003472236868877018),
  @(9, 3, 1, 3.154770000000001, 9.464310000000001, 0.03658830666704076, 0.04624159466596638, 3, 1, 1.525218666666667, 4.575656, 0.147299991145562, 0.1666904917928634, 4.811714093040001, 43.30542683736001, 0.00538945724808621, 0.007708034156156183),
  @(10, 3, 1, 3.154770000000001, 9.464310000000001, 0.03658830666704076, 0.04624159466596638, 3, 1, 1.729477333333333, 5.188432, 0.1670265395080728, 0.189013833582295, 5.456103206880001, 49.10492886192, 0.006111218249055966, 0.008740301078772907),
  @(11, 3, 1, 3.154770000000001, 9.464310000000001, 0.03658830666704076, 0.04624159466596638, 3, 1, 2.911407000000001, 8.734221000000002, 0.2811729456854671, 0.3181864182791614, 9.184819461390003, 82.66337515251003, 0.01028764196321506, 0.01471344738228062),
  @(12, 3, 1, 3.154770000000001, 9.464310000000001, 0.03658830666704076, 0.04624159466596638, 3, 1, 0.5748966666666666, 1.72469, 0.05552139883960665, 0.06283020932741304, 1.8136667571, 16.3230008139, 0.002031433967326609, 0.002905369072496054),
  @(13, 3, 1, 3.154770000000001, 9.464310000000001, 0.03658830666704076, 0.04624159466596638, 2, 1, 2.582909, 5.165818, 0.2494478209221534, 0.1881900088058249, 8.148483825930001, 48.89090295558, 0.009126873369324814, 0.008702206107383601),
  @(14, 3, 1, 2.963441333333333, 8.890324, 0.03436932020203823, 0.04343716117256439, 2, 1, 1.0305975, 2.061195, 0.09953130389913815, 0.07508903821244231, 3.05411522953, 18.32469137718, 0.003420823253835855, 0.003261654655126703),
  @(15, 3, 1, 2.963441333333333, 8.890324, 0.03436932020203823, 0.04343716117256439, 3, 1, 1.525218666666667, 4.575656, 0.147299991145562, 0.1666904917928634, 4.519896039171556, 40.679064352544, 0.005062600561439215, 0.007240561757940627),
  @(16, 3, 1, 2.963441333333333, 8.890324, 0.03436932020203823, 0.04343716117256439, 3, 1, 1.729477333333333, 5.188432, 0.1670265395080728, 0.189013833582295, 5.125204614663111, 46.126841531968, 0.005740588618591342, 0.008210224353158409),
  @(17, 3, 1, 2.963441333333333, 8.890324, 0.03436932020203823, 0.04343716117256439, 3, 1, 2.911407000000001, 8.734221000000002, 0.2811729456854671, 0.3181864182791614, 8.627783841956001, 77.65005457760401, 0.009663723002414122, 0.01382111473371292),
  @(18, 3, 1, 2.963441333333333, 8.890324, 0.03436932020203823, 0.04343716117256439, 3, 1, 0.5748966666666666, 1.72469, 0.05552139883960665, 0.06283020932741304, 1.703672544395555, 15.33305289956, 0.001908232734783515, 0.002729165929060798),
  @(19, 3, 1, 2.963441333333333, 8.890324, 0.03436932020203823, 0.04343716117256439, 2, 1, 2.582909, 5.165818, 0.2494478209221534, 0.1881900088058249, 7.654299290838667, 45.925795745032, 0.008573352030974183, 0.008174439743564929),
  @(20, 3, 1, 7.153525333333334, 21.460576, 0.08296496373632466, 0.1048540524021472, 2, 1, 1.0305975, 2.061195, 0.09953130389913815, 0.07508903821244231, 7.372405324720001, 44.23443194832, 0.008257611018621106, 0.007873389947554262),
  @(21, 3, 1, 7.153525333333334, 21.460576, 0.08296496373632466, 0.1048540524021472, 3, 1, 1.525218666666667, 4.575656, 0.147299991145562, 0.1666904917928634, 10.91069037087289, 98.19621333785601, 0.01222073842375249, 0.01747817356138859),
  @(22, 3, 1, 7.153525333333334, 21.460576, 0.08296496373632466, 0.1048540524021472, 3, 1, 1.729477333333333, 5.188432, 0.1670265395080728, 0.189013833582295, 12.37185991742578, 111.346739256832, 0.01385735079329106, 0.01981886641116869),
  @(23, 3, 1, 7.153525333333334, 21.460576, 0.08296496373632466, 0.1048540524021472, 3, 1, 2.911407000000001, 8.734221000000002, 0.2811729456854671, 0.3181864182791614, 20.826823730144, 187.441413571296, 0.02332750324243036, 0.03336313537589473),
  @(24, 3, 1, 7.153525333333334, 21.460576, 0.08296496373632466, 0.1048540524021472, 3, 1, 0.5748966666666666, 1.72469, 0.05552139883960665, 0.06283020932741304, 4.112537869048889, 37.01284082143999, 0.004606330841317984, 0.006588002061254447),
  @(25, 3, 1, 7.153525333333334, 21.460576, 0.08296496373632466, 0.1048540524021472, 2, 1, 2.582909, 5.165818, 0.2494478209221534, 0.1881900088058249, 18.47690496519467, 110.861429791168, 0.02069542941691167, 0.01973248504488652),
  @(26, 3, 1, 18.952291, 56.85687299999999, 0.2198043708894773, 0.2777965298305241, 2, 1, 1.0305975, 2.061195, 0.09953130389913815, 0.07508903821244231, 19.5321837238725, 117.193102343235, 0.02187741563735944, 0.0208594742437281),
  @(27, 3, 1, 18.952291, 56.85687299999999, 0.2198043708894773, 0.2777965298305241, 3, 1, 1.525218666666667, 4.575656, 0.147299991145562, 0.1666904917928634, 28.90638800929867, 260.157492083688, 0.03237718188577583, 0.0463060401758009),
  @(28, 3, 1, 18.952291, 56.85687299999999, 0.2198043708894773, 0.2777965298305241, 3, 1, 1.729477333333333, 5.188432, 0.1670265395080728, 0.189013833582295, 32.77755769923733, 294.9980192931359, 0.03671316343841836, 0.05250738705912572),
  @(29, 3, 1, 18.952291, 56.85687299999999, 0.2198043708894773, 0.2777965298305241, 3, 1, 2.911407000000001, 8.734221000000002, 0.2811729456854671, 0.3181864182791614, 55.177832683437, 496.600494150933, 0.06180304243753527, 0.08839108283715469),
  @(30, 3, 1, 18.952291, 56.85687299999999, 0.2198043708894773, 0.2777965298305241, 3, 1, 0.5748966666666666, 1.72469, 0.05552139883960665, 0.06283020932741304, 10.89560892159666, 98.06048029436998, 0.0122038461428435, 0.01745401411968077),
  @(31, 3, 1, 18.952291, 56.85687299999999, 0.2198043708894773, 0.2777965298305241, 2, 1, 2.582909, 5.165818, 0.2494478209221534, 0.1881900088058249, 48.952042994519, 293.7122579671139, 0.05482972134754494, 0.05227853139503394),
  @(32, 2, 1, 15.6101945, 31.220389, 0.1810434939783733, 0.1525394427540724, 2, 1, 1.0305975, 2.061195, 0.09953130389913815, 0.07508903821244231, 16.08782742621375, 64.351309704855, 0.01801949501812326, 0.0114540400458652),
  @(33, 2, 1, 15.6101945, 31.220389, 0.1810434939783733, 0.1525394427540724, 3, 1, 1.525218666666667, 4.575656, 0.147299991145562, 0.1666904917928634, 23.80896004169733, 142.853760250184, 0.02666770505997599, 0.02542687473048566),
  @(34, 2, 1, 15.6101945, 31.220389, 0.1810434939783733, 0.1525394427540724, 3, 1, 1.729477333333333, 5.188432, 0.1670265395080728, 0.189013833582295, 26.99747755667466, 161.984865340048, 0.03023906829965831, 0.02883206484745426),
  @(35, 2, 1, 15.6101945, 31.220389, 0.1810434939783733, 0.1525394427540724, 3, 1, 2.911407000000001, 8.734221000000002, 0.2811729456854671, 0.3181864182791614, 45.4476295386615, 272.685777231969, 0.05090453249908835, 0.04853597893621749),
  @(36, 2, 1, 15.6101945, 31.220389, 0.1810434939783733, 0.1525394427540724, 3, 1, 0.5748966666666666, 1.72469, 0.05552139883960665, 0.06283020932741304, 8.974248784068331, 53.84549270440999, 0.01005178803648919, 0.00958408511892531),
  @(37, 2, 1, 15.6101945, 31.220389, 0.1810434939783733, 0.1525394427540724, 2, 1, 2.582909, 5.165818, 0.2494478209221534, 0.1881900088058249, 40.31971186580049, 161.278847463202, 0.04516090506503823, 0.02870639907512452)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $c = 5
    for ($i = 1; $i -lt $entry.Length; $i++) {
        $ws.Cells.Item($r, $c).Value = $entry[$i]
        $c = $c + 1
    }
}
